# Insert a new weekly price record as row 696 ("Pintón", 5-Apr-2023 / serial 45021),
# pushing the existing rows 696-736 down to 697-737 (dimension grows to A1:T737).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 696 (shifts everything below it down by one).
$ws.Rows.Item(696).Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(696, 1).Value = 11
$ws.Cells.Item(696, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(696, 3).Value = "Bíobío"
$ws.Cells.Item(696, 4).Value = 45021
$ws.Cells.Item(696, 5).Value = 8
$ws.Cells.Item(696, 6).Value = "Fruta"
$ws.Cells.Item(696, 7).Value = 100108
$ws.Cells.Item(696, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(696, 9).Value = 100108006
$ws.Cells.Item(696, 10).Value = "Plátano"
$ws.Cells.Item(696, 11).Value = "Sin especificar"
$ws.Cells.Item(696, 12).Value = "Pintón"
$ws.Cells.Item(696, 13).Value = 1050
$ws.Cells.Item(696, 14).Value = 21000
$ws.Cells.Item(696, 15).Value = 22000
$ws.Cells.Item(696, 16).Value = 21524
$ws.Cells.Item(696, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(696, 18).Value = "Ecuador"
$ws.Cells.Item(696, 19).Value = 1076
$ws.Cells.Item(696, 20).Value = 20
